$p = $ppt.ActivePresentation

# Add a new slide after the existing one, using the "Title and Content"
# layout (PpSlideLayout.ppLayoutText = 2), matching slideLayout2.xml.
$p.Slides.Add(2, 2) | Out-Null
